# Daten aktualisiert am 2023-12-07
# Append three new ticker rows to the bottom of column A on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newValues = @("TAO-USD", "IMX-USD", "GRT-USD")

# Find the last used row in column A and append right after it.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($ws.Cells.Item(1, 1).Value -eq $null) {
    $lastRow = 0
}

foreach ($val in $newValues) {
    $lastRow = $lastRow + 1
    $ws.Cells.Item($lastRow, 1).Value = $val
}
